$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (A1:D1) to snake_case English names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case the Spanish connector words (de/del/la/el/los/las/y) within state
#    and municipality names, e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga"
$ws.Range("B9").Value = "Pabellón De Arteaga"
$ws.Range("B10").Value = "Rincón De Romos"
$ws.Range("B11").Value = "San Francisco De Los Romo"
$ws.Range("B12").Value = "San José De Gracia"
$ws.Range("B18").Value = "Playas De Rosarito"
$ws.Range("B46").Value = "Amatenango De La Frontera"
$ws.Range("B47").Value = "Amatenango Del Valle"
$ws.Range("B51").Value = "Bejucal De Ocampo"
$ws.Range("B53").Value = "Benemérito De Las Américas"
$ws.Range("B63").Value = "Chiapa De Corzo"
$ws.Range("B70").Value = "Comitán De Domínguez"
$ws.Range("B100").Value = "Marqués De Comillas"
$ws.Range("B101").Value = "Mazapa De Madero"
$ws.Range("B105").Value = "Montecristo De Guerrero"
$ws.Range("B109").Value = "Ocozocoautla De Espinosa"
$ws.Range("B122").Value = "Salto De Agua"
$ws.Range("B123").Value = "San Cristóbal De Las Casas"
$ws.Range("B181").Value = "Guadalupe Y Calvo"
$ws.Range("B185").Value = "Hidalgo Del Parral"
$ws.Range("B209").Value = "San Francisco De Borja"
$ws.Range("B210").Value = "San Francisco Del Oro"
$ws.Range("B217").Value = "Valle De Zaragoza"
$ws.Range("A219").Value = "Ciudad De México"
$ws.Range("B222").Value = "Cuajimalpa De Morelos"
$ws.Range("A237").Value = "Coahuila De Zaragoza"
$ws.Range("B252").Value = "San Juan De Sabinas"
$ws.Range("B268").Value = "Villa De Álvarez"
$ws.Range("B272").Value = "Coneto De Comonfort"
$ws.Range("B287").Value = "Nombre De Dios"
$ws.Range("B294").Value = "Pánuco De Coronado"
$ws.Range("B298").Value = "San Juan De Guadalupe"
$ws.Range("B299").Value = "San Juan Del Río"
$ws.Range("B300").Value = "San Luis Del Cordero"
$ws.Range("B301").Value = "San Pedro Del Gallo"
$ws.Range("A311").Value = "Estado De México"
$ws.Range("B311").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B314").Value = "Almoloya De Alquisiras"
$ws.Range("B315").Value = "Almoloya De Juárez"
$ws.Range("B316").Value = "Almoloya Del Río"
$ws.Range("B323").Value = "Atizapán De Zaragoza"
$ws.Range("B331").Value = "Chapa De Mota"
$ws.Range("B337").Value = "Coacalco De Berriozábal"
$ws.Range("B345").Value = "Ecatepec De Morelos"
$ws.Range("B353").Value = "Ixtapan De La Sal"
$ws.Range("B354").Value = "Ixtapan Del Oro"
$ws.Range("B371").Value = "Naucalpan De Juárez"
$ws.Range("B385").Value = "San Antonio La Isla"
$ws.Range("B386").Value = "San Felipe Del Progreso"
$ws.Range("B387").Value = "San José Del Rincón"
$ws.Range("B388").Value = "San Martín De Las Pirámides"
$ws.Range("B390").Value = "San Simón De Guerrero"
$ws.Range("B392").Value = "Soyaniquilpan De Juárez"
$ws.Range("B402").Value = "Tenango Del Aire"
$ws.Range("B403").Value = "Tenango Del Valle"
$ws.Range("B417").Value = "Tlalnepantla De Baz"
$ws.Range("B424").Value = "Valle De Bravo"
$ws.Range("B425").Value = "Valle De Chalco Solidaridad"
$ws.Range("B428").Value = "Villa De Allende"
$ws.Range("B429").Value = "Villa Del Carbón"
$ws.Range("B440").Value = "Apaseo El Alto"
$ws.Range("B441").Value = "Apaseo El Grande"
$ws.Range("B450").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B454").Value = "Jaral Del Progreso"
$ws.Range("B461").Value = "Purísima Del Rincón"
$ws.Range("B466").Value = "San Diego De La Unión"
$ws.Range("B468").Value = "San Francisco Del Rincón"
$ws.Range("B470").Value = "San Luis De La Paz"
$ws.Range("B471").Value = "San Miguel De Allende"
$ws.Range("B473").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B475").Value = "Silao De La Victoria"
$ws.Range("B480").Value = "Valle De Santiago"
$ws.Range("B486").Value = "Acapulco De Juárez"
$ws.Range("B489").Value = "Ajuchitlán Del Progreso"
$ws.Range("B490").Value = "Alcozauca De Guerrero"
$ws.Range("B494").Value = "Atenango Del Río"
$ws.Range("B495").Value = "Atlamajalcingo Del Monte"
$ws.Range("B497").Value = "Atoyac De Álvarez"
$ws.Range("B498").Value = "Ayutla De Los Libres"
$ws.Range("B501").Value = "Buenavista De Cuéllar"
$ws.Range("B502").Value = "Chilapa De Álvarez"
$ws.Range("B503").Value = "Chilpancingo De Los Bravo"
$ws.Range("B504").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B505").Value = "Cochoapa El Grande"
$ws.Range("B510").Value = "Coyuca De Benítez"
$ws.Range("B511").Value = "Coyuca De Catalán"
$ws.Range("B515").Value = "Cuetzala Del Progreso"
$ws.Range("B516").Value = "Cutzamala De Pinzón"
$ws.Range("B523").Value = "Huitzuco De Los Figueroa"
$ws.Range("B524").Value = "Iguala De La Independencia"
$ws.Range("B527").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B528").Value = "José Joaquín De Herrera"
$ws.Range("B531").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B537").Value = "Mártir De Cuilapan"
$ws.Range("B548").Value = "Taxco De Alarcón"
$ws.Range("B551").Value = "Tepecoacuilco De Trujano"
$ws.Range("B553").Value = "Tixtla De Guerrero"
$ws.Range("B557").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B558").Value = "Tlapa De Comonfort"
$ws.Range("B560").Value = "Técpan De Galeana"
$ws.Range("B565").Value = "Zihuatanejo De Azueta"
$ws.Range("B572").Value = "Agua Blanca De Iturbide"
$ws.Range("B579").Value = "Atotonilco De Tula"
$ws.Range("B580").Value = "Atotonilco El Grande"
$ws.Range("B586").Value = "Cuautepec De Hinojosa"
$ws.Range("B593").Value = "Huasca De Ocampo"
$ws.Range("B597").Value = "Huejutla De Reyes"
$ws.Range("B600").Value = "Jacala De Ledezma"
$ws.Range("B607").Value = "Mineral De La Reforma"
$ws.Range("B608").Value = "Mineral Del Chico"
$ws.Range("B609").Value = "Mineral Del Monte"
$ws.Range("B610").Value = "Mixquiahuala De Juárez"
$ws.Range("B611").Value = "Molango De Escamilla"
$ws.Range("B613").Value = "Nopala De Villagrán"
$ws.Range("B614").Value = "Omitlán De Juárez"
$ws.Range("B615").Value = "Pachuca De Soto"
$ws.Range("B618").Value = "Progreso De Obregón"
$ws.Range("B624").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B625").Value = "Santiago De Anaya"
$ws.Range("B629").Value = "Tenango De Doria"
$ws.Range("B631").Value = "Tepehuacán De Guerrero"
$ws.Range("B632").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B635").Value = "Tezontepec De Aldama"
$ws.Range("B644").Value = "Tula De Allende"
$ws.Range("B645").Value = "Tulancingo De Bravo"
$ws.Range("B646").Value = "Villa De Tezontepec"
$ws.Range("B650").Value = "Zacualtipán De Ángeles"
$ws.Range("B651").Value = "Zapotlán De Juárez"
$ws.Range("B656").Value = "Acatlán De Juárez"
$ws.Range("B657").Value = "Ahualulco De Mercado"
$ws.Range("B662").Value = "Atemajac De Brizuela"
$ws.Range("B665").Value = "Atotonilco El Alto"
$ws.Range("B667").Value = "Autlán De Navarro"
$ws.Range("B673").Value = "Cañadas De Obregón"
$ws.Range("B680").Value = "Concepción De Buenos Aires"
$ws.Range("B681").Value = "Cuautitlán De García Barragán"
$ws.Range("B691").Value = "Encarnación De Díaz"
$ws.Range("B697").Value = "Huejuquilla El Alto"
$ws.Range("B699").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B700").Value = "Ixtlahuacán Del Río"
$ws.Range("B704").Value = "Jilotlán De Los Dolores"
$ws.Range("B710").Value = "La Manzanilla De La Paz"
$ws.Range("B711").Value = "Lagos De Moreno"
$ws.Range("B719").Value = "Ojuelos De Jalisco"
$ws.Range("B724").Value = "San Cristóbal De La Barranca"
$ws.Range("B725").Value = "San Diego De Alejandría"
$ws.Range("B728").Value = "San Juan De Los Lagos"
$ws.Range("B729").Value = "San Juanito De Escobedo"
$ws.Range("B733").Value = "San Martín De Bolaños"
$ws.Range("B734").Value = "San Miguel El Alto"
$ws.Range("B736").Value = "San Sebastián Del Oeste"
$ws.Range("B737").Value = "Santa María De Los Ángeles"
$ws.Range("B738").Value = "Santa María Del Oro"
$ws.Range("B741").Value = "Talpa De Allende"
$ws.Range("B742").Value = "Tamazula De Gordiano"
$ws.Range("B745").Value = "Techaluta De Montenegro"
$ws.Range("B749").Value = "Teocuitatlán De Corona"
$ws.Range("B750").Value = "Tepatitlán De Morelos"
$ws.Range("B753").Value = "Tizapán El Alto"
$ws.Range("B754").Value = "Tlajomulco De Zúñiga"
$ws.Range("B765").Value = "Unión De San Antonio"
$ws.Range("B766").Value = "Unión De Tula"
$ws.Range("B767").Value = "Valle De Guadalupe"
$ws.Range("B768").Value = "Valle De Juárez"
$ws.Range("B773").Value = "Yahualica De González Gallo"
$ws.Range("B774").Value = "Zacoalco De Torres"
$ws.Range("B777").Value = "Zapotitlán De Vadillo"
$ws.Range("B779").Value = "Zapotlán Del Rey"
$ws.Range("B780").Value = "Zapotlán El Grande"
$ws.Range("A782").Value = "Michoacán De Ocampo"
$ws.Range("B804").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B806").Value = "Cojumatlán De Régules"
$ws.Range("B875").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B903").Value = "Coatlán Del Río"
$ws.Range("B913").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B917").Value = "Puente De Ixtla"
$ws.Range("B923").Value = "Tetela Del Volcán"
$ws.Range("B925").Value = "Tlaltizapán De Zapata"
$ws.Range("B935").Value = "Zacualpan De Amilpas"
$ws.Range("B939").Value = "Amatlán De Cañas"
$ws.Range("B940").Value = "Bahía De Banderas"
$ws.Range("B945").Value = "Ixtlán Del Río"
$ws.Range("B952").Value = "Santa María Del Oro"
$ws.Range("B980").Value = "Mier Y Noriega"
$ws.Range("B986").Value = "San Nicolás De Los Garza"
$ws.Range("B991").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B999").Value = "Ayoquezco De Aldama"
$ws.Range("B1003").Value = "Capulálpam De Méndez"
$ws.Range("B1005").Value = "Chalcatongo De Hidalgo"
$ws.Range("B1006").Value = "Chiquihuitlán De Benito Juárez"
$ws.Range("B1008").Value = "Ciénega De Zimatlán"
$ws.Range("B1010").Value = "Coicoyán De Las Flores"
$ws.Range("B1013").Value = "Constancia Del Rosario"
$ws.Range("B1016").Value = "Cuilápam De Guerrero"
$ws.Range("B1017").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B1019").Value = "El Barrio De La Soledad"
$ws.Range("B1021").Value = "Eloxochitlán De Flores Magón"
$ws.Range("B1022").Value = "Fresnillo De Trujano"
$ws.Range("B1024").Value = "Guadalupe De Ramírez"
$ws.Range("B1025").Value = "Guelatao De Juárez"
$ws.Range("B1026").Value = "Guevea De Humboldt"
$ws.Range("B1027").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B1028").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B1029").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B1030").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B1031").Value = "Huajuapan De León"
$ws.Range("B1032").Value = "Huautla De Jiménez"
$ws.Range("B1034").Value = "Ixtlán De Juárez"
$ws.Range("B1047").Value = "Magdalena Yodocono De Porfirio Díaz"
$ws.Range("B1049").Value = "Mariscala De Juárez"
$ws.Range("B1051").Value = "Mazatlán Villa De Flores"
$ws.Range("B1053").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B1054").Value = "Mixistlán De La Reforma"
$ws.Range("B1056").Value = "Mártires De Tacubaya"
$ws.Range("B1059").Value = "Nejapa De Madero"
$ws.Range("B1061").Value = "Oaxaca De Juárez"
$ws.Range("B1062").Value = "Ocotlán De Morelos"
$ws.Range("B1063").Value = "Pinotepa De Don Luis"
$ws.Range("B1065").Value = "Putla Villa De Guerrero"
$ws.Range("B1066").Value = "Reforma De Pineda"
$ws.Range("B1068").Value = "Rojas De Cuauhtémoc"
$ws.Range("B1077").Value = "San Agustín De Las Juntas"
$ws.Range("B1095").Value = "San Antonino El Alto"
$ws.Range("B1101").Value = "San Antonio De La Cal"
$ws.Range("B1104").Value = "San Baltazar Yatzachi El Bajo"
$ws.Range("B1121").Value = "San Dionisio Del Mar"
$ws.Range("B1123").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B1140").Value = "San Francisco Del Mar"
$ws.Range("B1160").Value = "San José Del Peñasco"
$ws.Range("B1161").Value = "San José Del Progreso"
$ws.Range("B1169").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B1210").Value = "San Juan De Los Cués"
$ws.Range("B1211").Value = "San Juan Del Estado"
$ws.Range("B1212").Value = "San Juan Del Río"
$ws.Range("B1233").Value = "San Martín De Los Cansecos"
$ws.Range("B1242").Value = "San Mateo Del Mar"
$ws.Range("B1269").Value = "San Miguel Del Puerto"
$ws.Range("B1270").Value = "San Miguel Del Río"
$ws.Range("B1271").Value = "San Miguel El Grande"
$ws.Range("B1281").Value = "San Pablo Villa De Mitla"
$ws.Range("B1317").Value = "San Pedro El Alto"
$ws.Range("B1318").Value = "San Pedro Y San Pablo Ayutla"
$ws.Range("B1319").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B1320").Value = "San Pedro Y San Pablo Tequixtepec"
$ws.Range("B1342").Value = "Santa Ana Del Valle"
$ws.Range("B1361").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B1367").Value = "Santa Cruz De Bravo"
$ws.Range("B1370").Value = "Santa Inés De Zaragoza"
$ws.Range("B1371").Value = "Santa Inés Del Monte"
$ws.Range("B1375").Value = "Santa Lucía Del Camino"
$ws.Range("B1394").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B1425").Value = "Santa María Del Rosario"
$ws.Range("B1426").Value = "Santa María Del Tule"
$ws.Range("B1427").Value = "Santa María La Asunción"
$ws.Range("B1478").Value = "Santiago Del Río"
$ws.Range("B1499").Value = "Santo Domingo De Morelos"
$ws.Range("B1509").Value = "Sitio De Xitlapehua"
$ws.Range("B1511").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B1512").Value = "Tanetze De Zaragoza"
$ws.Range("B1514").Value = "Tataltepec De Valdés"
$ws.Range("B1515").Value = "Teococuilco De Marcos Pérez"
$ws.Range("B1516").Value = "Teotitlán De Flores Magón"
$ws.Range("B1517").Value = "Teotitlán Del Valle"
$ws.Range("B1519").Value = "Tepelmeme Villa De Morelos"
$ws.Range("B1520").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B1521").Value = "Tlacolula De Matamoros"
$ws.Range("B1523").Value = "Tlalixtac De Cabrera"
$ws.Range("B1524").Value = "Totontepec Villa De Morelos"
$ws.Range("B1530").Value = "Villa Sola De Vega"
$ws.Range("B1531").Value = "Villa Talea De Castro"
$ws.Range("B1532").Value = "Villa Tejúpam De La Unión"
$ws.Range("B1533").Value = "Villa De Chilapa De Díaz"
$ws.Range("B1534").Value = "Villa De Etla"
$ws.Range("B1535").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B1536").Value = "Villa De Tututepec"
$ws.Range("B1537").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B1538").Value = "Villa De Zaachila"
$ws.Range("B1541").Value = "Yutanduchi De Guerrero"
$ws.Range("B1544").Value = "Zimatlán De Álvarez"
$ws.Range("B1572").Value = "Ayotoxco De Guerrero"
$ws.Range("B1576").Value = "Chalchicomula De Sesma"
$ws.Range("B1586").Value = "Chila De La Sal"
$ws.Range("B1597").Value = "Cuapiaxtla De Madero"
$ws.Range("B1601").Value = "Cuayuca De Andrade"
$ws.Range("B1602").Value = "Cuetzalan Del Progreso"
$ws.Range("B1619").Value = "Huehuetlán El Chico"
$ws.Range("B1620").Value = "Huehuetlán El Grande"
$ws.Range("B1625").Value = "Huitzilan De Serdán"
$ws.Range("B1627").Value = "Ixcamilpa De Guerrero"
$ws.Range("B1631").Value = "Izúcar De Matamoros"
$ws.Range("B1642").Value = "Los Reyes De Juárez"
$ws.Range("B1643").Value = "Mazapiltepec De Juárez"
$ws.Range("B1656").Value = "Palmar De Bravo"
$ws.Range("B1666").Value = "San Diego La Mesa Tochimiltzingo"
$ws.Range("B1683").Value = "San Nicolás De Los Ranchos"
$ws.Range("B1688").Value = "San Salvador El Seco"
$ws.Range("B1689").Value = "San Salvador El Verde"
$ws.Range("B1697").Value = "Tecali De Herrera"
$ws.Range("B1705").Value = "Tepanco De López"
$ws.Range("B1706").Value = "Tepango De Rodríguez"
$ws.Range("B1707").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B1712").Value = "Tepexi De Rodríguez"
$ws.Range("B1714").Value = "Tepeyahualco De Cuauhtémoc"
$ws.Range("B1715").Value = "Tetela De Ocampo"
$ws.Range("B1716").Value = "Teteles De Avila Castillo"
$ws.Range("B1721").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1733").Value = "Totoltepec De Guerrero"
$ws.Range("B1735").Value = "Tuzamapan De Galeana"
$ws.Range("B1739").Value = "Xayacatlán De Bravo"
$ws.Range("B1746").Value = "Xochitlán De Vicente Suárez"
$ws.Range("B1753").Value = "Zapotitlán De Méndez"
$ws.Range("B1762").Value = "Amealco De Bonfil"
$ws.Range("B1764").Value = "Cadereyta De Montes"
$ws.Range("B1771").Value = "Jalpan De Serra"
$ws.Range("B1772").Value = "Landa De Matamoros"
$ws.Range("B1775").Value = "Pinal De Amoles"
$ws.Range("B1778").Value = "San Juan Del Río"
$ws.Range("B1794").Value = "Armadillo De Los Infante"
$ws.Range("B1795").Value = "Axtla De Terrazas"
$ws.Range("B1802").Value = "Ciudad Del Maíz"
$ws.Range("B1811").Value = "Mexquitic De Carmona"
$ws.Range("B1817").Value = "San Ciro De Acosta"
$ws.Range("B1822").Value = "Santa María Del Río"
$ws.Range("B1824").Value = "Soledad De Graciano Sánchez"
$ws.Range("B1832").Value = "Tanquián De Escobedo"
$ws.Range("B1837").Value = "Villa De Arista"
$ws.Range("B1838").Value = "Villa De Arriaga"
$ws.Range("B1839").Value = "Villa De Guadalupe"
$ws.Range("B1840").Value = "Villa De Ramos"
$ws.Range("B1841").Value = "Villa De Reyes"
$ws.Range("B1842").Value = "Villa De La Paz"
$ws.Range("B1900").Value = "Nacozari De García"
$ws.Range("B1913").Value = "San Miguel De Horcasitas"
$ws.Range("B1914").Value = "San Pedro De La Cueva"
$ws.Range("B1937").Value = "Jalpa De Méndez"
$ws.Range("B1972").Value = "Soto La Marina"
$ws.Range("B1979").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B1980").Value = "Amaxac De Guerrero"
$ws.Range("B1981").Value = "Apetatitlán De Antonio Carvajal"
$ws.Range("B1987").Value = "Contla De Juan Cuamatzi"
$ws.Range("B1996").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B2000").Value = "Muñoz De Domingo Arenas"
$ws.Range("B2001").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B2004").Value = "Papalotla De Xicohténcatl"
$ws.Range("B2008").Value = "San Pablo Del Monte"
$ws.Range("B2009").Value = "Sanctórum De Lázaro Cárdenas"
$ws.Range("B2017").Value = "Tepetitla De Lardizábal"
$ws.Range("B2020").Value = "Tetla De La Solidaridad"
$ws.Range("B2032").Value = "Ziltlaltépec De Trinidad Sánchez Santos"
$ws.Range("A2034").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B2042").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B2046").Value = "Amatlán De Los Reyes"
$ws.Range("B2058").Value = "Boca Del Río"
$ws.Range("B2060").Value = "Camarón De Tejeda"
$ws.Range("B2064").Value = "Castillo De Teayo"
$ws.Range("B2066").Value = "Cazones De Herrera"
$ws.Range("B2074").Value = "Chinampa De Gorostiza"
$ws.Range("B2086").Value = "Cosamaloapan De Carpio"
$ws.Range("B2087").Value = "Cosautlán De Carvajal"
$ws.Range("B2106").Value = "Hueyapan De Ocampo"
$ws.Range("B2107").Value = "Huiloapan De Cuauhtémoc"
$ws.Range("B2108").Value = "Ignacio De La Llave"
$ws.Range("B2112").Value = "Ixhuacán De Los Reyes"
$ws.Range("B2114").Value = "Ixhuatlán De Madero"
$ws.Range("B2115").Value = "Ixhuatlán Del Café"
$ws.Range("B2116").Value = "Ixhuatlán Del Sureste"
$ws.Range("B2126").Value = "Juchique De Ferrer"
$ws.Range("B2132").Value = "Las Vigas De Ramírez"
$ws.Range("B2133").Value = "Lerdo De Tejada"
$ws.Range("B2139").Value = "Martínez De La Torre"
$ws.Range("B2143").Value = "Medellín De Bravo"
$ws.Range("B2147").Value = "Mixtla De Altamirano"
$ws.Range("B2149").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B2160").Value = "Ozuluama De Mascareñas"
$ws.Range("B2163").Value = "Paso De Ovejas"
$ws.Range("B2164").Value = "Paso Del Macho"
$ws.Range("B2168").Value = "Poza Rica De Hidalgo"
$ws.Range("B2181").Value = "Sayula De Alemán"
$ws.Range("B2185").Value = "Soledad De Doblado"
$ws.Range("B2193").Value = "Tatahuicapan De Juárez"
$ws.Range("B2227").Value = "Vega De Alatorre"
$ws.Range("B2239").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B2240").Value = "Zozocolco De Hidalgo"
$ws.Range("B2264").Value = "Dzilam De Bravo"
$ws.Range("B2330").Value = "Cañitas De Felipe Pescador"
$ws.Range("B2332").Value = "Concepción Del Oro"
$ws.Range("B2335").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B2345").Value = "Jiménez Del Teul"
$ws.Range("B2351").Value = "Mezquital Del Oro"
$ws.Range("B2356").Value = "Moyahua De Estrada"
$ws.Range("B2357").Value = "Nochistlán De Mejía"
$ws.Range("B2358").Value = "Noria De Ángeles"
$ws.Range("B2364").Value = "Santa María De La Paz"
$ws.Range("B2370").Value = "Teúl De González Ortega"
$ws.Range("B2371").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B2373").Value = "Trinidad García De La Cadena"
$ws.Range("B2379").Value = "Villa De Cos"

# 3. Rename the grand-total label row from "TOTAL" to "Total"
$ws.Range("A2383").Value = "Total"

# 4. Remove the trailing footnote/metadata rows (sample size, source, etc.)
$ws.Rows("2385:2389").Delete()

